$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 290, shifting existing rows 290-315 down to 291-316.
$ws.Rows.Item(290).Insert()

# Fill the fixed (unchanged-across-the-series) columns in the newly inserted row.
$ws.Range("A290").Value = 10
$ws.Range("B290").Value = "Vega Modelo de Temuco"
$ws.Range("C290").Value = "La Araucanía"
$ws.Range("E290").Value = 9
$ws.Range("F290").Value = 100112001
$ws.Range("G290").Value = "Berenjena"
$ws.Range("H290").Value = "Sin especificar"
$ws.Range("I290").Value = "Primera"
$ws.Range("N290").Value = "`$/caja 60 unidades"
$ws.Range("Q290").Value = 60
$ws.Range("R290").Value = "Hortaliza"

# Match the date style/number format of column D from the neighboring row.
$ws.Range("D290").NumberFormat = $ws.Range("D291").NumberFormat

# New data values for the inserted row.
$ws.Range("D290").Value = 44769
$ws.Range("J290").Value = 300
$ws.Range("K290").Value = 13000
$ws.Range("L290").Value = 15000
$ws.Range("M290").Value = 14333
$ws.Range("O290").Value = "Región de Arica y Parinacota"
$ws.Range("P290").Value = 239
